$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 (CasesTab query): drop the trailing Cohort line/column from the RETURN clause
$caseQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN [''UBC01''] and demo.breed in [''Basset Hound'',''Belgian Malinois'', ''Labrador Retriever'',''West Highland White Terrier'']and diag.disease_term in [''Bladder Cancer''] and diag.primary_disease_site in [ ''Bladder, Prostate'', ''Bladder, Urethra'', ''Bladder, Urethra, Prostate'', ''Urethra, Prostate''] and diag.best_response in [''Not Determined'', ''Partial Response'',''Progressive Disease'',''Stable Disease'']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`'
$ws.Range("B2").Value = $caseQuery

# Row 2 shrinks now that the query text has one fewer line
$ws.Rows.Item(2).RowHeight = 333.5

# Move the active selection/top-left viewport to B2 (was D4/topLeftCell B4)
$ws.Range("B2").Select()

Write-Output "done"
